# Apply the "list-of-questions" update:
#  - Move the two-sample-t-test interpretation question data from row 16
#    (matched-pairs t-test row) down to row 13 (two-sample t-test row),
#    adding a new "...interpret-p" file to the D-column reference.
#  - Clear out the old row 16 cells (C16:E16) that used to hold that data.
#  - Rename the "mult choice" solution type to "schoice" on the two
#    linear-model rows (21 and 24).
#  - Update the active selection to C14 (matches the author's last edit
#    location on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move data onto row 13 (two-sample t-test) ---
$ws.Range("C13").Value = "read test-statistic, read sided t-test"
$ws.Range("D13").Value = "schoice-2samtt-interpret-t, schoice-2samtt-interpret-alt, schoice-2samtt-interpret-p"
$ws.Range("E13").Value = "single choice"

# --- Clear out the old row 16 (matched pairs t-test) cells ---
$ws.Range("C16:E16").ClearContents()

# --- Rename "mult choice" -> "schoice" for the lm rows ---
$ws.Range("E21").Value = "schoice"
$ws.Range("E24").Value = "schoice"

# --- Update the saved selection ---
$ws.Range("C14").Select()
